$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 914.1948
$ws.Range("J17").Value = 918.71234
$ws.Range("L17").Value = 2756.13702
$ws.Range("N17").Value = -3092.13702
$ws.Range("H33").Value = 244.4
$ws.Range("I33").Value = 235.92308
$ws.Range("K33").Value = 235.92308
$ws.Range("M33").Value = -6.923079999999999
$ws.Range("H40").Value = 4702.4736
$ws.Range("J40").Value = 2706
$ws.Range("L40").Value = 2706
$ws.Range("N40").Value = -3056
$ws.Range("H43").Value = 2528.7693
$ws.Range("I43").Value = 2894.25
$ws.Range("J43").Value = 2366.3333
$ws.Range("K43").Value = 2894.25
$ws.Range("L43").Value = 2366.3333
$ws.Range("M43").Value = -2825.25
$ws.Range("N43").Value = -2504.3333
$ws.Range("H88").Value = 7996.4287
$ws.Range("I88").Value = 4500
$ws.Range("J88").Value = 8579.166999999999
$ws.Range("K88").Value = 4500
$ws.Range("L88").Value = 8579.166999999999
$ws.Range("M88").Value = -4094
$ws.Range("N88").Value = -9391.166999999999
$ws.Range("H91").Value = 7996.4287
$ws.Range("I91").Value = 4500
$ws.Range("J91").Value = 8579.166999999999
$ws.Range("K91").Value = 4500
$ws.Range("L91").Value = 8579.166999999999
$ws.Range("M91").Value = -3096
$ws.Range("N91").Value = -11387.167
$ws.Range("H116").Value = 695783.8
$ws.Range("J116").Value = 1201817.6
$ws.Range("L116").Value = 1201817.6
$ws.Range("N116").Value = -1208701.6
$ws.Range("H125").Value = 8468.666999999999
$ws.Range("I125").Value = 2495
$ws.Range("K125").Value = 22455
$ws.Range("M125").Value = -19995
$ws.Range("H129").Value = 2279.5625
$ws.Range("J129").Value = 2750
$ws.Range("L129").Value = 8250
$ws.Range("N129").Value = -18250
$ws.Range("H132").Value = 48644.06
$ws.Range("I132").Value = 55581.836
$ws.Range("K132").Value = 166745.508
$ws.Range("M132").Value = -164215.508
$ws.Range("H135").Value = 3000.4
$ws.Range("I135").Value = 1000.5
$ws.Range("K135").Value = 9004.5
$ws.Range("M135").Value = -6469.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 55699.223
$ws.Range("J43").Value = 42762
$ws.Range("L43").Value = 42762
$ws.Range("N43").Value = -43388
$ws.Range("H74").Value = 5210578
$ws.Range("I74").Value = 6945935.5
$ws.Range("K74").Value = 6945935.5
$ws.Range("M74").Value = -6945061.5
$ws.Range("H77").Value = 5210578
$ws.Range("I77").Value = 6945935.5
$ws.Range("K77").Value = 34729677.5
$ws.Range("M77").Value = -34725309.5
$ws.Range("H97").Value = 753.34784
$ws.Range("I97").Value = 696.7273
$ws.Range("J97").Value = 1999
$ws.Range("K97").Value = 696.7273
$ws.Range("L97").Value = 1999
$ws.Range("M97").Value = -200.7273
$ws.Range("N97").Value = -2991
$ws.Range("H102").Value = 25156.934
$ws.Range("I102").Value = 26853.072
$ws.Range("K102").Value = 26853.072
$ws.Range("M102").Value = -25231.072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2513.1428
$ws.Range("I86").Value = 2499.3333
$ws.Range("J86").Value = 2523.5
$ws.Range("K86").Value = 2499.3333
$ws.Range("L86").Value = 2523.5
$ws.Range("M86").Value = -1376.3333
$ws.Range("N86").Value = -4769.5
$ws.Range("H89").Value = 2513.1428
$ws.Range("I89").Value = 2499.3333
$ws.Range("J89").Value = 2523.5
$ws.Range("K89").Value = 12496.6665
$ws.Range("L89").Value = 12617.5
$ws.Range("M89").Value = -6880.666499999999
$ws.Range("N89").Value = -23849.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 590979.3
$ws.Range("I58").Value = 727768.6
$ws.Range("K58").Value = 727768.6
$ws.Range("M58").Value = -727565.6
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 20000
$ws.Range("K60").Value = 20000
$ws.Range("M60").Value = -19489
$ws.Range("H99").Value = 2853
$ws.Range("I99").Value = 3316.125
$ws.Range("J99").Value = 2235.5
$ws.Range("K99").Value = 3316.125
$ws.Range("L99").Value = 2235.5
$ws.Range("M99").Value = -1818.125
$ws.Range("N99").Value = -5231.5
$ws.Range("H122").Value = 2054.5557
$ws.Range("J122").Value = 3284.8572
$ws.Range("L122").Value = 9854.571599999999
$ws.Range("N122").Value = -14754.5716
$ws.Range("H126").Value = 2853
$ws.Range("I126").Value = 3316.125
$ws.Range("J126").Value = 2235.5
$ws.Range("K126").Value = 9948.375
$ws.Range("L126").Value = 6706.5
$ws.Range("M126").Value = -7478.375
$ws.Range("N126").Value = -11646.5
$ws.Range("H132").Value = 25521236
$ws.Range("I132").Value = 30316742
$ws.Range("J132").Value = 15630506
$ws.Range("K132").Value = 90950226
$ws.Range("L132").Value = 46891518
$ws.Range("M132").Value = -90947696
$ws.Range("N132").Value = -46896578
$ws.Range("H134").Value = 17504.762
$ws.Range("I134").Value = 18899.264
$ws.Range("K134").Value = 56697.792
$ws.Range("M134").Value = -54162.792
$ws.Range("H136").Value = 590979.3
$ws.Range("I136").Value = 727768.6
$ws.Range("K136").Value = 2183305.8
$ws.Range("M136").Value = -2180755.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 416.66666
$ws.Range("J23").Value = 525
$ws.Range("L23").Value = 1575
$ws.Range("N23").Value = -2045
$ws.Range("H86").Value = 433
$ws.Range("J86").Value = 500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 433
$ws.Range("J89").Value = 500
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -16356
$ws.Range("H129").Value = 995.25
$ws.Range("I129").Value = 699.2727
$ws.Range("J129").Value = 1646.4
$ws.Range("K129").Value = 2097.8181
$ws.Range("L129").Value = 4939.200000000001
$ws.Range("M129").Value = 2902.1819
$ws.Range("N129").Value = -14939.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 750
$ws.Range("J25").Value = 750
$ws.Range("L25").Value = 750
$ws.Range("N25").Value = -1808
$ws.Range("H122").Value = 86847
$ws.Range("I122").Value = 170668.17
$ws.Range("J122").Value = 15000.286
$ws.Range("K122").Value = 512004.51
$ws.Range("L122").Value = 45000.858
$ws.Range("M122").Value = -509554.51
$ws.Range("N122").Value = -49900.858
$ws.Range("H126").Value = 796544.8
$ws.Range("I126").Value = 1043915.06
$ws.Range("K126").Value = 3131745.18
$ws.Range("M126").Value = -3129275.18

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3190.862
$ws.Range("I7").Value = 3061.4
$ws.Range("K7").Value = 3061.4
$ws.Range("M7").Value = -2949.4
$ws.Range("H40").Value = 5055.6924
$ws.Range("I40").Value = 5072.5
$ws.Range("J40").Value = 4999.6665
$ws.Range("K40").Value = 5072.5
$ws.Range("L40").Value = 4999.6665
$ws.Range("M40").Value = -4936.5
$ws.Range("N40").Value = -5271.6665
$ws.Range("H55").Value = 328.5
$ws.Range("I55").Value = 328.5
$ws.Range("K55").Value = 328.5
$ws.Range("M55").Value = -155.5
$ws.Range("H126").Value = 3190.862
$ws.Range("I126").Value = 3061.4
$ws.Range("K126").Value = 9184.200000000001
$ws.Range("M126").Value = -6714.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1168325.4
$ws.Range("I136").Value = 1386636.9
$ws.Range("K136").Value = 4159910.7
$ws.Range("M136").Value = -4157360.7
